$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "p0001"
$ws.Range("B2").Value = "pilas doble aa"
$ws.Range("C2").Value = 1565
$ws.Range("D2").Value = 94
$ws.Range("E2").Value = 56

$ws.Range("A3").Value = "p0002"
$ws.Range("B3").Value = "porta retrato 20*25"
$ws.Range("C3").Value = 6590
$ws.Range("D3").Value = 99
$ws.Range("E3").Value = 0
